$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old Category column (C) data rows; column C header stays.
$ws.Range("C2:C18").Clear()

# Rewrite String_ID / Translation rows (2-26) with the revised window action strings.
$ws.Range("A2").Value = "About_LicenceLabel"
$ws.Range("B2").Value = "License"
$ws.Range("A3").Value = "About_LicenceTab"
$ws.Range("B3").Value = "_License"
$ws.Range("A4").Value = "Action_MinimiseToolTip"
$ws.Range("B4").Value = "Minimize the window if it is already active"
$ws.Range("A5").Value = "Action_RestoreIfMinimised"
$ws.Range("B5").Value = "_Restore if minimized"
$ws.Range("A6").Value = "Custom_ButtonBgColorLabel"
$ws.Range("B6").Value = "Background color:"
$ws.Range("A7").Value = "Custom_ButtonBorderColourLabel"
$ws.Range("B7").Value = "Border color:"
$ws.Range("A8").Value = "Custom_ButtonTextColourLabel"
$ws.Range("B8").Value = "Text color:"
$ws.Range("A9").Value = "Custom_WindowBgColorLabel"
$ws.Range("B9").Value = "Background color:"
$ws.Range("A10").Value = "E_CUST007"
$ws.Range("B10").Value = "Error while selecting background color"
$ws.Range("A11").Value = "E_CUST017"
$ws.Range("B11").Value = "Error while changing button border color"
$ws.Range("A12").Value = "E_CUST018"
$ws.Range("B12").Value = "Error while changing button background color"
$ws.Range("A13").Value = "E_CUST024"
$ws.Range("B13").Value = "Error while changing button text color"
$ws.Range("A14").Value = "E_MAIN001"
$ws.Range("B14").Value = "Error while initializing folders"
$ws.Range("A15").Value = "E_MAIN002"
$ws.Range("B15").Value = "Error while initializing application"
$ws.Range("A16").Value = "E_REG015"
$ws.Range("B16").Value = "Error while changing the region color"
$ws.Range("A17").Value = "Options_ModeOverlayColorLabel"
$ws.Range("B17").Value = "Mode/page overlay text color:"
$ws.Range("A18").Value = "Options_PointerColourLabel"
$ws.Range("B18").Value = "Pointer indicator color:"
$ws.Range("A19").Value = "Regions_BackgroundColorLabel"
$ws.Range("B19").Value = "Background color:"
$ws.Range("A20").Value = "Regions_RegionColorLabel"
$ws.Range("B20").Value = "Color:"
$ws.Range("A21").Value = "String_Maximise"
$ws.Range("B21").Value = "Maximize"
$ws.Range("A22").Value = "String_MaximiseOrMinimiseWindow"
$ws.Range("B22").Value = "Maximize / minimize window"
$ws.Range("A23").Value = "String_MaximiseOrRestoreWindow"
$ws.Range("B23").Value = "Maximize / restore window"
$ws.Range("A24").Value = "String_MaximiseWindow"
$ws.Range("B24").Value = "Maximize window"
$ws.Range("A25").Value = "String_Minimise"
$ws.Range("B25").Value = "Minimize"
$ws.Range("A26").Value = "String_MinimiseWindow"
$ws.Range("B26").Value = "Minimize window"

# Move the active selection to A2 (matches the saved view state).
[void]$ws.Range("A2").Select()
